$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing row 2 data (will be rewritten) ---
$ws.Range("A2:M2").ClearContents()

# --- Header row (row 1) ---
$ws.Range("A1").Value = 'ID'
$ws.Range("B1").Value = 'ID'
$ws.Range("C1").Value = 'Product Name'
$ws.Range("D1").Value = 'Category'
$ws.Range("E1").Value = 'Main Image'
$ws.Range("F1").Value = 'Description'
$ws.Range("G1").Value = 'Sowing Time'
$ws.Range("I1").Value = 'Fertilizer Info'
$ws.Range("J1").Value = 'Fertilizer Info'

# --- Clear header cells that no longer hold text (H1, K1, L1) but keep bold style ---
$ws.Range("H1").ClearContents()
$ws.Range("K1").ClearContents()
$ws.Range("L1").ClearContents()
$ws.Range("M1").ClearContents()

# --- Row 2: Bhindo (ભીંડો) ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 'ભીંડો '
$ws.Range("C2").Value = 'vegetable'
$ws.Range("D2").Value = 'assets/images/products/1.png'
$ws.Range("E2").Value = 'ભીંડો સરિતા ખેતી વિષયક માહિતી '
$ws.Range("F2").Value = '<div style="background-color: #fff; border-radius: 8px; box-shadow: 0 2px 4px rgba(0,0,0,0.1); padding: 20px; width: 300px; font-family: ''Roboto'', sans-serif; margin: 20px auto;">
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવેતર સમય:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> ઠંડીના દિવસો સિવાય બારેમાસ.</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવણી અંતર:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> બે હાર (લાઈન) વચ્ચેનું અંતર ૧૮ થી ૨૪ ઇંચ તથા બે છોડ વચ્ચેનું અંતર ૬ થી ૯ ઇંચ<br>
    બીજનું પ્રમાણ પ્રતિ એકરે ૨ કિલો થી ૨ કિલો ૫૦૦ ગ્રામ</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div>
    <span style="font-weight: bold; font-size: 18px; color: #000;">જમીન:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> સામાન્ય રીતે દરેક પ્રકારની જમીન અનુકૂળ રહે છે.</span>
  </div>
</div>
'
$ws.Range("I2").Value = '<div style="font-family: ''Noto Sans Gujarati'', sans-serif; font-size: 20px; line-height: 1.9; color: #333; padding: 20px;">
  <p><strong style="color: #006400;">ખાતર:</strong> પાયાના ખાતર તરીકે છાણીયા ખાતર ની જરૂર રહે છે સાથે રાસાયણિક ખાતરોમાં પ્રતિ એકરે 
ડાયએમોનિયમ ફોસ્ફેટ (ડી.એ.પી.) ખાતર ૧૮ કિલો તથા ૧૨ કિલો પોટાશ ખાતરની જરૂરિયાત રહે છે. </p>
  <p>પછી ના રાસાયણિક ખાતરનો ઉપયોગ જમીન ચકાસણી વિશ્લેષણ પર આધારિત હોવો જોઈએ.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #006400;">પાક સંરક્ષણ:</strong> બીજ ઊગી નીકળ્યા બાદ નિંદામણની કાળજી રાખવી ચોળીના પાકમાં આવતા મુખ્ય રોગોમાં થીપ્સ, સફેદ ચૂસ્યા, લીલા ચૂસ્યા, કથીરી, ઈયળ તેમજ ફૂગજન્ય રોગો આવી શકે છે. કૃષિ તજજ્ઞ ની માહિતી લઈને જરૂરિયાત મુજબ પાક સંરક્ષણ પગલાં સમયસર લેવા.  </p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #006400;">ભીંડો ઉગ્યાથી:</strong> ૪૫ થી ૫૦ દિવસે છોડનો વધુ પડતો વિકાસ જણાઈ તો ખાસ આલ્ફા નેફથાઈલ એસિટીક એસિડ ૧૦ મિલી અથવા ક્લોરમેકવાટ ક્લોરાઇડ ૧૫ મિલી - ૧૫ લીટર પાણીમાં નાખી છટકાવ કરી શકાય તેમજ બીજો છટકાવ ૮ થી ૧૦ દિવસે કરી શકાય.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #006400;">પિયત:</strong> બીજ ઉગી નીકળ્યા બાદ જમીનનો પ્રકાર તથા તાપમાન મુજબ પાણી આપવું વરસાદના સમયમાં તેમજ નદીના પટમાં લેવાયેલ પાકને પિયત આપવાની જરૂરિયાત ઓછી રહે છે. </p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #b30000;">ખુલાસો:</strong> આ પત્રિકામાં આપેલ માર્ગદર્શનમાં જમીનનો પ્રકાર પ્રતિકૂળ આબોહવાની પરિસ્થિતિ અને ઋતુ અપર્યાપ્ત જંતુઓ અને રોગોના હુમલા પાક અને ઉત્પાદન પર પ્રતિકૂળ અસર પેદા કરી શકે છે પાક નિયંત્રણ અમારા નિયંત્રણ બહાર છે તેથી ખેડૂત પોતે ઉત્પાદન ઉપજ માટે સંપૂર્ણપણે જવાબદાર છે તેમ છતાં અમે તેમને યોગ્ય સમયે યોગ્ય બીજ અને શ્રેષ્ઠ પાકના પાકને વાવણી દ્વારા મહત્તમ ઉપજ મેળવવા માટે સૂચવીએ છીએ માર્ગદર્શિકાઓના પગલે પણ કંપની પર પાક નિષ્ફળતા માટે કોઈ જવાબદાર નથી</p>
</div>
'

# --- Row 3: Guvar (ગુવાર) ---
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 'ગુવાર '
$ws.Range("C3").Value = 'vegetable'
$ws.Range("D3").Value = 'assets/images/products/Guwar S 11.png'
$ws.Range("E3").Value = 'ગુવાર S.11 ખેતી વિશે માહિતી'
$ws.Range("F3").Value = '<div style="background-color: #fff; border-radius: 8px; box-shadow: 0 2px 4px rgba(0,0,0,0.1); padding: 20px; width: 300px; font-family: ''Roboto'', sans-serif; margin: 20px auto;">
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવેતર સમય:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> ઠંડીના દિવસો સિવાય બારેમાસ</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવણી અંતર:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> બે હાર લાઈન વચ્ચેનું અંતર ૧૮ ઇંચ તથા બે છોડ વચ્ચેનું અંતર ૩ થી ૪ ઇંચ.</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">બીજ નું પ્રમાણ:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> પ્રતિ એકરે ૫ થી ૬ કિલો.</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div>
    <span style="font-weight: bold; font-size: 18px; color: #000;">જમીન:</span>
    <span style="font-size: 16px; color: #3e7e3b;"> સામાન્ય રીતે દરેક પ્રકારની જમીન અનુકૂળ આવે છે તેમજ જમીન ઊંડી ખેડ તથા પોચી ભરભરી કરવી અનુકૂળ નથી.</span>
  </div>
</div>
'
$ws.Range("I3").Value = '<div style="font-family: ''Noto Sans Gujarati'', sans-serif; font-size: 20px; line-height: 1.9; color: #333; padding: 20px;">
  <p><strong style="color: #006400;">ખાતર:</strong> પાયાના ખાતર તરીકે છાણિયા ખાતર ની જરૂર રહે છે સાથે રાસાયણિક ખાતરોનું પ્રતિએકરે ડાયએમોનિયમ ફોસ્ફેટ (ડીએપી) ખાતર ની જરૂરિયાત રહે છે પછીના રાસાયણિક ખાતરનો ઉપયોગ જમીન ચકાસણી વિશ્લેષણ પર આધારિત હોવો જોઈએ.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #006400;">પાક સંરક્ષણ:</strong> બીજ ઊગી નીકળ્યા બાદ નિંદામણ કાળજી રાખવી ગુવારના પાકમાં આવતા મુખ્ય રોગોમાં સફેદ ચૂસ્યા, લીલા ચૂસ્યા તેમજ વાયરસ જન્ય, ફૂગ જન્ય રોગો આવી શકે છે કૃષિ તજજ્ઞ ની માહિતી લઈને જરૂરિયાત મુજબ પાક સંરક્ષણ પગલા સમયસર લેવા.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #006400;">પિયત:</strong> બીજ ઉગી નીકળ્યા બાદ જમીનનો પ્રકાર તથા તાપમાન મુજબ પાણી આપવું વરસાદના સમયમાં તેમજ નદીના પટમાં લેવાયેલ પાકને પિયત આપવાની જરૂરિયાત ઓછી રહે છે.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #b30000;">ખુલાસો:</strong> આ પત્રિકામાં આપેલ માર્ગદર્શનમાં જમીનનો પ્રકાર પ્રતિકૂળ આબોહવાની પરિસ્થિતિ અને ઋતુ અપર્યાપ્ત જંતુઓ અને રોગોના હુમલા પાક અને ઉત્પાદન પર પ્રતિકૂળ અસર પેદા કરી શકે છે પાક નિયંત્રણ અમારા નિયંત્રણ બહાર છે તેથી ખેડૂત પોતે ઉત્પાદન ઉપજ માટે સંપૂર્ણપણે જવાબદાર છે તેમ છતાં અમે તેમને યોગ્ય સમયે યોગ્ય બીજ અને શ્રેષ્ઠ પાકના પાકને વાવણી દ્વારા મહત્તમ ઉપજ મેળવવા માટે સૂચવીએ છીએ માર્ગદર્શિકાઓના પગલે પણ કંપની પર પાક નિષ્ફળતા માટે કોઈ જવાબદાર નથી.</p>
</div>
'

# --- Row 4: Choli (ચોળી) ---
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 'ચોળી '
$ws.Range("C4").Value = 'vegetable'
$ws.Range("D4").Value = 'assets/images/products/CHOLIE SOMNATH.png'
$ws.Range("E4").Value = 'ચોળી સોમનાથ ખેતી વિષયક માહિતી'
$ws.Range("F4").Value = '<div style="background-color: #fff; border-radius: 8px; box-shadow: 0 2px 4px rgba(0,0,0,0.1); padding: 20px; width: 90%; max-width: 800px; font-family: ''Roboto'', sans-serif; margin: 20px auto;">
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">ફળનો રંગ તથા આકાર:</span>
    <span style="font-size: 16px; color: #3e7e3b;">સિંગ પાતળી આકર્ષક, લીલો કલર તેમજ ૫ થી ૬ ઇંચ લાંબી</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવેતર સમય:</span>
    <span style="font-size: 16px; color: #3e7e3b;">ફેબ્રુઆરી, માર્ચ, એપ્રિલ, ઓગસ્ટ, સપ્ટેમ્બર તથા નવેમ્બર</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવણી અંતર:</span>
    <span style="font-size: 16px; color: #3e7e3b;">બે હાર (લાઈન) વચ્ચે ૧૮ થી ૨૦ ઇંચ અને બે છોડ વચ્ચે ૬ થી ૭ ઇંચ</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવેતર સમય:</span>
    <span style="font-size: 16px; color: #3e7e3b;">મે, જૂન, જુલાઈ</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">વાવણી અંતર:</span>
    <span style="font-size: 16px; color: #3e7e3b;">બે હાર લાઈન વચ્ચે ૨૦ થી ૨૪ ઇંચ તથા બે છોડ વચ્ચે ૮ થી ૯ ઇંચ</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">બીજનું પ્રમાણ:</span>
    <span style="font-size: 16px; color: #3e7e3b;">પ્રતિ એકરે ૬ થી ૭ કિલો</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div style="margin-bottom: 20px;">
    <span style="font-weight: bold; font-size: 18px; color: #000;">પહેલી વીણી:</span>
    <span style="font-size: 16px; color: #3e7e3b;">૫૫ થી ૬૦ દિવસમાં</span>
  </div>
  <div style="border-top: 1px solid #eee; margin: 10px 0;"></div>
  <div>
    <span style="font-weight: bold; font-size: 18px; color: #000;">જમીન:</span>
    <span style="font-size: 16px; color: #3e7e3b;">સામાન્ય રીતે દરેક પ્રકારની જમીન અનુકૂળ રહે છે</span>
  </div>
</div>
'
$ws.Range("I4").Value = '<div style="font-family: ''Noto Sans Gujarati'', sans-serif; font-size: 20px; line-height: 1.9; color: #333; padding: 20px;">
  <p><strong style="color: #006400;">ખાતર:</strong> પાયાના ખાતર તરીકે છાણીયા ખાતર ની જરૂર રહે છે સાથે રાસાયણિક ખાતરોમાં પ્રતિ એકરે ડાયએમોનિયમ ફોસ્ફેટ (ડીએપી) ખાતર ૧૮ કિલો તથા ૧૨ કિલો પોટાશ ખાતરની જરૂરિયાત રહે છે. પછી ના રાસાયણિક ખાતરનો ઉપયોગ જમીન ચકાસણી વિશ્લેષણ પર આધારિત હોવો જોઈએ.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #006400;">પાક સંરક્ષણ:</strong> બીજ ઊગી નીકળ્યા બાદ નિંદામણની કાળજી રાખવી ચોળીના પાકમાં આવતા મુખ્ય રોગોમાં થીપ્સ, સફેદ ચૂસ્યા, લીલા ચૂસ્યા, કથીરી ઈયળ તેમજ ફૂગજન્ય રોગો આવી શકે છે. કૃષિ તજજ્ઞ ની માહિતી લઈને જરૂરિયાત મુજબ પાક સંરક્ષણ પગલાં સમયસર લેવા.</p>
  <p><strong style="color: #006400;">ચોળી ઉગ્યાથી:</strong> ૪૫ થી ૫૦ દિવસે છોડનો વધુ પડતો વિકાસ જણાઈ તો ખાસ આલ્ફા નેફથાઈલ એસિટીક એસિડ ૧૦ મિલી અથવા ક્લોરમેકવાટ ક્લોરાઇડ ૧૫ મિલી - ૧૫ લીટર પાણીમાં નાખી છટકાવ કરી શકાય તેમજ બીજો છટકાવ ૮ થી ૧૦ દિવસે કરી શકાય.</p>
  <p><strong style="color: #006400;">પિયત:</strong> બીજ ઉગી નીકળ્યા બાદ જમીનનો પ્રકાર તથા તાપમાન મુજબ પાણી આપવું વરસાદના સમયમાં તેમજ નદીના પટમાં લેવાયેલ પાકને પિયત આપવાની જરૂરિયાત ઓછી રહે છે. ખાસ ચોળી સોમનાથ મિશ્ર પાક અન્ય પાક સાથે વાવેલ ન હોય તો જ.</p>
  <hr style="border: 1px solid #006400; margin: 20px 0;">
  <p><strong style="color: #b30000;">ખુલાસો:</strong> આ પત્રિકામાં આપેલ માર્ગદર્શનમાં જમીનનો પ્રકાર પ્રતિકૂળ આબોહવાની પરિસ્થિતિ અને ઋતુ અપર્યાપ્ત જંતુઓ અને રોગોના હુમલા પાક અને ઉત્પાદન પર પ્રતિકૂળ અસર પેદા કરી શકે છે પાક નિયંત્રણ અમારા નિયંત્રણ બહાર છે તેથી ખેડૂત પોતે ઉત્પાદન ઉપજ માટે સંપૂર્ણપણે જવાબદાર છે તેમ છતાં અમે તેમને યોગ્ય સમયે યોગ્ય બીજ અને શ્રેષ્ઠ પાકના પાકને વાવણી દ્વારા મહત્તમ ઉપજ મેળવવા માટે સૂચવીએ છીએ માર્ગદર્શિકાઓના પગલે પણ કંપની પર પાક નિષ્ફળતા માટે કોઈ જવાબદાર નથી.</p>
</div>
'

# --- Reset auto-expanded row heights caused by embedded newlines in the HTML blobs ---
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# --- Column widths (COM ColumnWidth snaps to pixel grid; these inputs land on the
#     closest achievable stored widths to the target 24 / 25.88671875) ---
$ws.Columns.Item(9).ColumnWidth = 23.09
$ws.Columns.Item(10).ColumnWidth = 25

# --- Selection ---
$ws.Range("I1").Select()
